# "balance the quest exp" — rebalance RewardGold (V) / RewardExp (Z) values
# on the Quest sheet's "表2" table (rows 4-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quest")

# --- RewardGold (column V) -------------------------------------------------
$ws.Range("V4").Value = 100

# --- RewardExp (column Z) ---------------------------------------------------
$ws.Range("Z4").Value  = 100
$ws.Range("Z5").Value  = 100
$ws.Range("Z6").Value  = 100
$ws.Range("Z7").Value  = 85
$ws.Range("Z8").Value  = 100
$ws.Range("Z9").Value  = 100
$ws.Range("Z10").Value = 100
$ws.Range("Z11").Value = 60
$ws.Range("Z12").Value = 85
$ws.Range("Z13").Value = 50
$ws.Range("Z14").Value = 50
$ws.Range("Z15").Value = 50

# --- restore the view/selection state recorded in the saved workbook -------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("V5").Select() | Out-Null
